# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the 66182380-... file has
# just been re-handed-off for both zh-cn and de-de, so its status flips from
# "In Translation" to "Ready for handoff", its priority flips from "ht" to
# "mt", and the handoff timestamps + the Overview roll-up date move forward.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"              # zh-cn column
$overview.Range("F3").Value = "Ready for handoff"              # de-de column
$overview.Range("G3").Value = "2016-09-02 20:16:37"            # Latest HO Xliff Generate Date

# ---- zh-cn sheet ----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"                  # Status
$zhcn.Range("E3").Value = "mt"                                 # Priority
$zhcn.Range("H3").Value = "2016-09-02 20:16:30"                # Latest Handoff Datetime
$zhcn.Columns.Item(3).ColumnWidth = 16.38                       # Status column widened to fit new text

# ---- de-de sheet ----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"                   # Status
$dede.Range("E3").Value = "mt"                                  # Priority
$dede.Range("H3").Value = "2016-09-02 20:16:37"                 # Latest Handoff Datetime
$dede.Columns.Item(3).ColumnWidth = 16.38                        # Status column widened to fit new text
